{"js": "const mapping = [\n  [\"2025-11-16 Sunday\", \"2025-11-17 Monday\"],\n  [\"147\u00f76=24, 3\", \"521\u00f78=65, 1\"],\n  [\"914\u00f75=182, 4\", \"535\u00f73=178, 1\"],\n  [\"873\u00f77=124, 5\", \"836\u00f73=278, 2\"],\n  [\"153\u00f72=76, 1\", \"186\u00f75=37, 1\"],\n  [\"562\u00f72=281, 0\", \"346\u00f74=86, 2\"],\n  [\"235\u00f74=58, 3\", \"802\u00f79=89, 1\"],\n  [\"795\u00f73=265, 0\", \"216\u00f79=24, 0\"],\n  [\"503\u00f72=251, 1\", \"776\u00f76=129, 2\"],\n  [\"891\u00f77=127, 2\", \"410\u00f74=102, 2\"],\n  [\"267\u00f76=44, 3\", \"474\u00f79=52, 6\"],\n  [\"568\u00f77=81, 1\", \"138\u00f72=69, 0\"],\n  [\"676\u00f73=225, 1\", \"767\u00f78=95, 7\"],\n  [\"405\u00f78=50, 5\", \"621\u00f72=310, 1\"],\n  [\"583\u00f75=116, 3\", \"310\u00f72=155, 0\"],\n  [\"869\u00f73=289, 2\", \"746\u00f79=82, 8\"],\n  [\"934\u00f72=467, 0\", \"641\u00f75=128, 1\"],\n  [\"542\u00f75=108, 2\", \"748\u00f79=83, 1\"],\n  [\"814\u00f74=203, 2\", \"940\u00f79=104, 4\"],\n  [\"763\u00f77=109, 0\", \"394\u00f77=56, 2\"],\n  [\"840\u00f79=93, 3\", \"688\u00f75=137, 3\"],\n  [\"239\u00f77=34, 1\", \"847\u00f79=94, 1\"],\n  [\"761\u00f72=380, 1\", \"816\u00f73=272, 0\"],\n  [\"988\u00f79=109, 7\", \"843\u00f76=140, 3\"],\n  [\"589\u00f73=196, 1\", \"710\u00f77=101, 3\"],\n  [\"972\u00f78=121, 4\", \"529\u00f77=75, 4\"],\n];\n\nfor (const [oldText, newText] of mapping) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old='2025-11-16 Sunday'; new='2025-11-17 Monday'},\n    @{old='147\u00f76=24, 3'; new='521\u00f78=65, 1'},\n    @{old='914\u00f75=182, 4'; new='535\u00f73=178, 1'},\n    @{old='873\u00f77=124, 5'; new='836\u00f73=278, 2'},\n    @{old='153\u00f72=76, 1'; new='186\u00f75=37, 1'},\n    @{old='562\u00f72=281, 0'; new='346\u00f74=86, 2'},\n    @{old='235\u00f74=58, 3'; new='802\u00f79=89, 1'},\n    @{old='795\u00f73=265, 0'; new='216\u00f79=24, 0'},\n    @{old='503\u00f72=251, 1'; new='776\u00f76=129, 2'},\n    @{old='891\u00f77=127, 2'; new='410\u00f74=102, 2'},\n    @{old='267\u00f76=44, 3'; new='474\u00f79=52, 6'},\n    @{old='568\u00f77=81, 1'; new='138\u00f72=69, 0'},\n    @{old='676\u00f73=225, 1'; new='767\u00f78=95, 7'},\n    @{old='405\u00f78=50, 5'; new='621\u00f72=310, 1'},\n    @{old='583\u00f75=116, 3'; new='310\u00f72=155, 0'},\n    @{old='869\u00f73=289, 2'; new='746\u00f79=82, 8'},\n    @{old='934\u00f72=467, 0'; new='641\u00f75=128, 1'},\n    @{old='542\u00f75=108, 2'; new='748\u00f79=83, 1'},\n    @{old='814\u00f74=203, 2'; new='940\u00f79=104, 4'},\n    @{old='763\u00f77=109, 0'; new='394\u00f77=56, 2'},\n    @{old='840\u00f79=93, 3'; new='688\u00f75=137, 3'},\n    @{old='239\u00f77=34, 1'; new='847\u00f79=94, 1'},\n    @{old='761\u00f72=380, 1'; new='816\u00f73=272, 0'},\n    @{old='988\u00f79=109, 7'; new='843\u00f76=140, 3'},\n    @{old='589\u00f73=196, 1'; new='710\u00f77=101, 3'},\n    @{old='972\u00f78=121, 4'; new='529\u00f77=75, 4'},\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.old, $false, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}"}
